$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ("Clockwise") updates ---
# Remove the red highlight fill from H4:M4 (matches the "no fill" look of B4:G4)
$ws.Range("H4:M4").Interior.Pattern = -4142

# Update counts / timings across H4:M4
$ws.Range("H4").Value = 5

$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "0.00683689117432 seconds"
$ws.Range("I4").NumberFormat = "General"

$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 0.00939416885376
$ws.Range("L4").Value = 5

$ws.Range("M4").NumberFormat = "@"
$ws.Range("M4").Value = "0.00760984420776"
$ws.Range("M4").NumberFormat = "General"

# --- Row 8 ("Hybrid") updates ---
$ws.Range("K8").Value = 5.9993429184
$ws.Range("M8").Value = 6.45539999008

# --- Selection moved to K9 ---
$ws.Range("K9").Select()
